# Adapt column header formatting to respective input file names.
# - Rename header cells in row 1 from "<Name>_old" / "<Name>_new" suffixes
#   to "<Name>_FV2310" / "<Name>_FV2404" respectively.
# - Turn the data range into an Excel Table ("Table1").
# - Freeze the header row (row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old header text (row 1, columns A:U) -> new header text
$oldHeaders = @(
    "Segmentname_old",
    "Segmentgruppe_old",
    "Segment_old",
    "Datenelement_old",
    "Segment ID_old",
    "Code_old",
    "Qualifier_old",
    "Beschreibung_old",
    "Bedingungsausdruck_old",
    "Bedingung_old",
    "diff",
    "Segmentname_new",
    "Segmentgruppe_new",
    "Segment_new",
    "Datenelement_new",
    "Segment ID_new",
    "Code_new",
    "Qualifier_new",
    "Beschreibung_new",
    "Bedingungsausdruck_new",
    "Bedingung_new"
)

$newHeaders = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

$headerRow = $ws.Range("A1:U1")
for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $oldText = $oldHeaders[$i]
    $newText = $newHeaders[$i]
    if ($oldText -ne $newText) {
        $headerRow.Replace($oldText, $newText, 1, 1, $false, $false, $false, $false) | Out-Null
    }
}

# Convert the used range into a native Excel table so the headers become
# table column names, and a filter/autofilter is shown.
$dataRange = $ws.UsedRange
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"

# Freeze the header row (split below row 1).
$ws.Activate() | Out-Null
$ws.Cells.Item(2, 1).Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
